# Update loading_percent values on Sheet1 for rows 2-25 (Case with 380 kV)
# Columns B, C, E, F, G, H, J, L, M are updated per the source diff;
# columns A, D, I, K, N, O are unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$BValues = @(18.3493160628503, 17.96304917179411, 17.72568667873684, 17.62905275092927, 17.61301659731319, 17.72438287266877, 18.21625525604764, 19.17325635128007, 19.86404489104534, 20.17411782276275, 20.2908184289458, 20.26571855874285, 20.18373387402394, 20.13341911455726, 19.8436867543706, 19.66479117669135, 19.561509213524, 19.52647713338362, 19.68387573751695, 20.20783512985922, 20.54604283006023, 20.36595895487548, 19.67524894571854, 18.91599064409192)
for ($i = 0; $i -lt $BValues.Length; $i++) {
    $ws.Cells.Item(2 + $i, "B").Value = $BValues[$i]
}

$CValues = @(10.13689606908455, 9.880661142585586, 9.718174030785473, 9.650715248277894, 9.639440261716075, 9.717269221890245, 10.04964736817116, 10.6583868451902, 11.07680972055649, 11.26043518074786, 11.32896896725866, 11.31425404361856, 11.26609375181894, 11.2364628097135, 11.06467110459962, 10.95753583343444, 10.8952840878159, 10.87409958051491, 10.9690060731948, 11.28026703308995, 11.47784437324162, 11.37293965258421, 10.96382242544753, 10.49859612496616)
for ($i = 0; $i -lt $CValues.Length; $i++) {
    $ws.Cells.Item(2 + $i, "C").Value = $CValues[$i]
}

$EValues = @(14.31650721457518, 14.33240089419212, 14.34285432992622, 14.34728925622496, 14.34803625747217, 14.34291343151795, 14.3218434143766, 14.28601958776148, 14.26302577508131, 14.2532825567825, 14.24969574695224, 14.25046366675829, 14.25298541071951, 14.25454342098532, 14.26367691189576, 14.26946335666365, 14.27285905346549, 14.27402037932673, 14.26884039755056, 14.25224192784159, 14.24199253638131, 14.24740816140854, 14.26912182245298, 14.29512512279172)
for ($i = 0; $i -lt $EValues.Length; $i++) {
    $ws.Cells.Item(2 + $i, "E").Value = $EValues[$i]
}

$FValues = @(42.59027401372364, 42.61331257678801, 42.63809905264106, 42.65087019932982, 42.65315194702333, 42.63826048427504, 42.59600649243534, 42.59775197119266, 42.65076764827555, 42.68612100642124, 42.7011214666341, 42.69781917718402, 42.68732280029599, 42.68110340257181, 42.64868311904029, 42.63167069010466, 42.62294347510861, 42.62017034272983, 42.63337222094963, 42.69036209747226, 42.73700708230136, 42.71125314725245, 42.63259967724425, 42.58820004141089)
for ($i = 0; $i -lt $FValues.Length; $i++) {
    $ws.Cells.Item(2 + $i, "F").Value = $FValues[$i]
}

$GValues = @(46.13213581801816, 46.06539803120612, 46.04042980811631, 46.03427594200512, 46.03349672071391, 46.04033054432315, 46.10580049838266, 46.3612102861285, 46.62599005964792, 46.76302799195534, 46.81728448045398, 46.8054946735353, 46.76744448748209, 46.74444467024903, 46.61736631818353, 46.54364223331224, 46.50280108716242, 46.48924202195691, 46.55132866987903, 46.77855681275854, 46.94082452726537, 46.85296858382433, 46.54784882146621, 46.27852910080739)
for ($i = 0; $i -lt $GValues.Length; $i++) {
    $ws.Cells.Item(2 + $i, "G").Value = $GValues[$i]
}

$HValues = @(18.46127899789551, 18.50677104642342, 18.53874886742334, 18.55279405459972, 18.55518738816279, 18.53893418478698, 18.47612306374394, 18.38519212644782, 18.33822563453939, 18.32120337644906, 18.31538459773125, 18.31660984186119, 18.32071207796183, 18.32330656853199, 18.33942570441257, 18.3504285893089, 18.35716579976113, 18.359517003427, 18.349215000708, 18.31949011019982, 18.30371993887341, 18.3118013626946, 18.34976238249708, 18.40631980453043)
for ($i = 0; $i -lt $HValues.Length; $i++) {
    $ws.Cells.Item(2 + $i, "H").Value = $HValues[$i]
}

$JValues = @(8.905232056508321, 8.924320961054518, 8.936566704456848, 8.941689534398749, 8.942548200809544, 8.936635255160521, 8.911705265710403, 8.866958372729947, 8.836571099030662, 8.823279731754239, 8.8183225457941, 8.819386794241671, 8.822870381620634, 8.825014056902818, 8.837450380238074, 8.845215531767234, 8.849731943190333, 8.851269743627675, 8.844383736668476, 8.821845110054841, 8.807557361430085, 8.815142682613587, 8.844759628826624, 8.878624033393315)
for ($i = 0; $i -lt $JValues.Length; $i++) {
    $ws.Cells.Item(2 + $i, "J").Value = $JValues[$i]
}

$LValues = @(11.81608460278643, 11.818061323912, 11.82051042864313, 11.82182014746824, 11.82205647982099, 11.82052682861475, 11.81651040906387, 11.81838976402655, 11.82565347372796, 11.83022051798358, 11.83213013370555, 11.83171087352829, 11.83037402052155, 11.82957858075372, 11.82538029759018, 11.82312722634403, 11.82195026323098, 11.82157223425529, 11.82335476937306, 11.83076180827769, 11.83665223768927, 11.83341285632113, 11.82325152840491, 11.81684334171432)
for ($i = 0; $i -lt $LValues.Length; $i++) {
    $ws.Cells.Item(2 + $i, "L").Value = $LValues[$i]
}

$MValues = @(17.47772239151322, 17.40533899344867, 17.363613606642, 17.34730594006519, 17.34464042931711, 17.36339084304939, 17.45220796812588, 17.64737748847064, 17.80276405482154, 17.87587321084625, 17.90388965895772, 17.89784134078404, 17.87817159330621, 17.86616598100259, 17.7980335295678, 17.75684525691672, 17.73338395937316, 17.72548022122622, 17.76120622750755, 17.88394021200464, 17.96607863592572, 17.92206938751781, 17.75923395117565, 17.59241140855192)
for ($i = 0; $i -lt $MValues.Length; $i++) {
    $ws.Cells.Item(2 + $i, "M").Value = $MValues[$i]
}
